$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.686896920204163
$ws.Range("B1").Value = 3.321837186813354
$ws.Range("C1").Value = 6.138805866241455
$ws.Range("D1").Value = 1.809961438179016
$ws.Range("E1").Value = 0.8947840332984924
